# Weekly data refresh: insert a new record as row 385 (most recent week),
# pushing the existing rows 385-427 down to 386-428.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(385).Insert()

$ws.Range("A385").Value = 10
$ws.Range("B385").Value = "Vega Modelo de Temuco"
$ws.Range("C385").Value = "La Araucanía"
$ws.Range("D385").Value = 44918
$ws.Range("E385").Value = 9
$ws.Range("F385").Value = 100112009
$ws.Range("G385").Value = "Acelga"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 90
$ws.Range("K385").Value = 9000
$ws.Range("L385").Value = 10000
$ws.Range("M385").Value = 9389
$ws.Range("N385").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O385").Value = "Provincia de Cautín"
$ws.Range("P385").Value = 782
$ws.Range("Q385").Value = 12
$ws.Range("R385").Value = "Hortaliza"
